$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3325, 4090, 4289, 4410, 4521, 4720, 5102, 5130, 5130, 5359, 5359, 5359, 5398, 5446)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
